$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data. Each cell is forced to remain plain
# text (matching the source feed inlineStr cells) by temporarily marking
# it Text-formatted before the write, then clearing the format override so
# the cell style reverts to the original default (no lingering at-format).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "71.997.74"
Set-TextValue "E2" "  +4.17%  "
Set-TextValue "D3" "4.039.67"
Set-TextValue "E3" "  +3.85%  "
Set-TextValue "E4" "  -0.07%  "
Set-TextValue "D5" "516.98"
Set-TextValue "E5" "  -1.71%  "
Set-TextValue "D6" "148.67"
Set-TextValue "E6" "  +4.18%  "
Set-TextValue "D7" "0.620"
Set-TextValue "E7" "  +2.22%  "
Set-TextValue "E8" "  +0.04%  "
Set-TextValue "D9" "0.735"
Set-TextValue "E9" "  +2.76%  "
Set-TextValue "E10" "  +2.73%  "
Set-TextValue "E11" "  +1.24%  "
Set-TextValue "D12" "47.98"
Set-TextValue "E12" "  +14.96%  "
Set-TextValue "D13" "10.82"
Set-TextValue "E13" "  +6.85%  "
Set-TextValue "D14" "4.687.88"
Set-TextValue "E14" "  +3.69%  "
Set-TextValue "D15" "4.073.78"
Set-TextValue "E15" "  +3.98%  "
Set-TextValue "D16" "21.17"
Set-TextValue "E16" "  +8.19%  "
Set-TextValue "D17" "14.10"
Set-TextValue "E17" "  +3.30%  "
Set-TextValue "E18" "  +0.01%  "
Set-TextValue "E19" "  -2.34%  "
Set-TextValue "D20" "72.006.87"
Set-TextValue "E20" "  +4.34%  "
Set-TextValue "D21" "435.78"
Set-TextValue "E21" "  +2.99%  "
Set-TextValue "D22" "96.40"
Set-TextValue "E22" "  +10.59%  "
Set-TextValue "D23" "3.51"
Set-TextValue "E23" "  +6.56%  "
Set-TextValue "D24" "14.56"
Set-TextValue "E24" "  +3.32%  "
Set-TextValue "B25" "RenderToken"
Set-TextValue "C25" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D25" "11.88"
Set-TextValue "E25" "  +4.38%  "
Set-TextValue "B26" "PancakeSwap"
Set-TextValue "C26" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D26" "3.99"
Set-TextValue "E26" "  -2.41%  "
Set-TextValue "D27" "11.15"
Set-TextValue "E27" "  +6.47%  "
Set-TextValue "D28" "36.84"
Set-TextValue "E28" "  +3.38%  "
Set-TextValue "E29" "  +10.80%  "
Set-TextValue "D30" "698.20"
Set-TextValue "E30" "  +0.69%  "
Set-TextValue "D31" "13.48"
Set-TextValue "E31" "  +3.76%  "
Set-TextValue "E32" "  +3.33%  "
Set-TextValue "D33" "6.99"
Set-TextValue "E33" "  +19.19%  "
Set-TextValue "D34" "67.94"
Set-TextValue "E34" "  -0.54%  "
Set-TextValue "D35" "0.0₃0892"
Set-TextValue "E35" "  +7.27%  "
Set-TextValue "B36" "ThetaToken"
Set-TextValue "C36" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D36" "3.67"
Set-TextValue "E36" "  +25.70%  "
Set-TextValue "B37" "TheGraph"
Set-TextValue "C37" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D37" "0.434"
Set-TextValue "E37" "  -1.95%  "
Set-TextValue "D38" "40.36"
Set-TextValue "E38" "  +1.69%  "
Set-TextValue "E39" "  +5.00%  "
Set-TextValue "D40" "0.999"
Set-TextValue "D41" "1.00"
Set-TextValue "E41" "  -0.09%  "
Set-TextValue "D42" "0.0485"
Set-TextValue "E42" "  +2.06%  "
Set-TextValue "E43" "  +3.36%  "
Set-TextValue "E44" "  +0.53%  "
Set-TextValue "E45" "  +5.60%  "
Set-TextValue "D46" "0.144"
Set-TextValue "E46" "  +3.97%  "
Set-TextValue "D47" "3.11"
Set-TextValue "E47" "  +3.12%  "
Set-TextValue "D48" "9.01"
Set-TextValue "E48" "  +8.85%  "
Set-TextValue "D49" "3.32"
Set-TextValue "E49" "  +2.21%  "
Set-TextValue "D50" "0.000269"
Set-TextValue "E50" "  +20.40%  "
Set-TextValue "D51" "0.0₆0338"
Set-TextValue "E51" "  +4.29%  "
